$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1424.5
$ws.Range("I33").Value = 728.8889
$ws.Range("J33").Value = 3511.3333
$ws.Range("K33").Value = 728.8889
$ws.Range("L33").Value = 3511.3333
$ws.Range("M33").Value = -499.8889
$ws.Range("N33").Value = -3969.3333
$ws.Range("H34").Value = 5850
$ws.Range("I34").Value = 1133.3334
$ws.Range("J34").Value = 20000
$ws.Range("K34").Value = 1133.3334
$ws.Range("L34").Value = 20000
$ws.Range("M34").Value = -930.3334
$ws.Range("N34").Value = -20406
$ws.Range("H36").Value = 5850
$ws.Range("I36").Value = 1133.3334
$ws.Range("J36").Value = 20000
$ws.Range("K36").Value = 1133.3334
$ws.Range("L36").Value = 20000
$ws.Range("M36").Value = -418.3334
$ws.Range("N36").Value = -21430
$ws.Range("H86").Value = 38465210
$ws.Range("I86").Value = 83335900
$ws.Range("J86").Value = 4614.2856
$ws.Range("K86").Value = 83335900
$ws.Range("L86").Value = 4614.2856
$ws.Range("M86").Value = -83334777
$ws.Range("N86").Value = -6860.2856
$ws.Range("H89").Value = 38465210
$ws.Range("I89").Value = 83335900
$ws.Range("J89").Value = 4614.2856
$ws.Range("K89").Value = 416679500
$ws.Range("L89").Value = 23071.428
$ws.Range("M89").Value = -416673884
$ws.Range("N89").Value = -34303.428
$ws.Range("H138").Value = 4632314
$ws.Range("I138").Value = 8131861.5
$ws.Range("J138").Value = 3880.1936
$ws.Range("K138").Value = 24395584.5
$ws.Range("L138").Value = 11640.5808
$ws.Range("M138").Value = -24390444.5
$ws.Range("N138").Value = -21920.5808
$ws.Range("H141").Value = 1471.2106
$ws.Range("I141").Value = 1409
$ws.Range("J141").Value = 2000
$ws.Range("K141").Value = 4227
$ws.Range("L141").Value = 6000
$ws.Range("M141").Value = 953
$ws.Range("N141").Value = -16360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1869.7142
$ws.Range("I2").Value = 1898.3334
$ws.Range("J2").Value = 1848.25
$ws.Range("K2").Value = 1898.3334
$ws.Range("L2").Value = 1848.25
$ws.Range("M2").Value = -1785.3334
$ws.Range("N2").Value = -2074.25
$ws.Range("H102").Value = 2676.5
$ws.Range("I102").Value = 2504
$ws.Range("J102").Value = 2849
$ws.Range("K102").Value = 2504
$ws.Range("L102").Value = 2849
$ws.Range("M102").Value = -882
$ws.Range("N102").Value = -6093
$ws.Range("H116").Value = 1869.7142
$ws.Range("I116").Value = 1898.3334
$ws.Range("J116").Value = 1848.25
$ws.Range("K116").Value = 1898.3334
$ws.Range("L116").Value = 1848.25
$ws.Range("M116").Value = 395.6666
$ws.Range("N116").Value = -6436.25
$ws.Range("H122").Value = 4636.6045
$ws.Range("I122").Value = 4904.757
$ws.Range("J122").Value = 2983
$ws.Range("K122").Value = 14714.271
$ws.Range("L122").Value = 8949
$ws.Range("M122").Value = -12264.271
$ws.Range("N122").Value = -13849

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1869.7142
$ws.Range("I3").Value = 1898.3334
$ws.Range("J3").Value = 1848.25
$ws.Range("K3").Value = 1898.3334
$ws.Range("L3").Value = 1848.25
$ws.Range("M3").Value = -1784.3334
$ws.Range("N3").Value = -2076.25
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
$ws.Range("H94").Value = 777.74286
$ws.Range("I94").Value = 582.8929000000001
$ws.Range("K94").Value = 582.8929000000001
$ws.Range("M94").Value = -131.8929000000001
$ws.Range("H140").Value = 41426.668
$ws.Range("J140").Value = 41426.668
$ws.Range("L140").Value = 41426.668
$ws.Range("N140").Value = -51786.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1608.2051
$ws.Range("I58").Value = 604.8
$ws.Range("J58").Value = 3400
$ws.Range("K58").Value = 604.8
$ws.Range("L58").Value = 3400
$ws.Range("M58").Value = -401.8
$ws.Range("N58").Value = -3806
$ws.Range("H64").Value = 24271
$ws.Range("J64").Value = 24271
$ws.Range("L64").Value = 24271
$ws.Range("N64").Value = -24767
$ws.Range("H67").Value = 24271
$ws.Range("J67").Value = 24271
$ws.Range("L67").Value = 24271
$ws.Range("N67").Value = -25987
$ws.Range("H134").Value = 459195
$ws.Range("I134").Value = 1343.9744
$ws.Range("J134").Value = 1832748.1
$ws.Range("K134").Value = 4031.9232
$ws.Range("L134").Value = 5498244.300000001
$ws.Range("M134").Value = -1496.9232
$ws.Range("N134").Value = -5503314.300000001
$ws.Range("H136").Value = 1608.2051
$ws.Range("I136").Value = 604.8
$ws.Range("J136").Value = 3400
$ws.Range("K136").Value = 1814.4
$ws.Range("L136").Value = 10200
$ws.Range("M136").Value = 735.6000000000001
$ws.Range("N136").Value = -15300

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 10000
$ws.Range("I57").Value = 2000
$ws.Range("J57").Value = 14000
$ws.Range("K57").Value = 6000
$ws.Range("L57").Value = 42000
$ws.Range("M57").Value = -5441
$ws.Range("N57").Value = -43118

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1506.8
$ws.Range("I97").Value = 1321.4286
$ws.Range("J97").Value = 2480
$ws.Range("K97").Value = 1321.4286
$ws.Range("L97").Value = 2480
$ws.Range("M97").Value = -825.4286
$ws.Range("N97").Value = -3472
$ws.Range("H126").Value = 3661.8948
$ws.Range("I126").Value = 2639.7144
$ws.Range("J126").Value = 4258.1665
$ws.Range("K126").Value = 7919.1432
$ws.Range("L126").Value = 12774.4995
$ws.Range("M126").Value = -5449.1432
$ws.Range("N126").Value = -17714.4995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3433.3333
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 3433.3333
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 6866.6666
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -8988.6666
$ws.Range("H84").Value = 3433.3333
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 3433.3333
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 34333.333
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -44941.333
